# Scheduled-runner market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on affected rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets with freshly polled prices, per Lamia_Profits.xlsx upstream diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1937.8667
$ws.Range("J17").Value = 1969.1428
$ws.Range("L17").Value = 5907.428400000001
$ws.Range("N17").Value = -6243.428400000001

# Row 19
$ws.Range("H19").Value = 1563.2273
$ws.Range("I19").Value = 2418.3333
$ws.Range("J19").Value = 537.1
$ws.Range("K19").Value = 2418.3333
$ws.Range("L19").Value = 537.1
$ws.Range("M19").Value = -2243.3333
$ws.Range("N19").Value = -887.1

# Row 41
$ws.Range("H41").Value = 225.5
$ws.Range("I41").Value = 98.75
$ws.Range("J41").Value = 310
$ws.Range("K41").Value = 98.75
$ws.Range("L41").Value = 310
$ws.Range("M41").Value = 341.25
$ws.Range("N41").Value = -1190

# Row 43
$ws.Range("H43").Value = 5861.8
$ws.Range("I43").Value = 744
$ws.Range("K43").Value = 744
$ws.Range("M43").Value = -675

# Row 62
$ws.Range("H62").Value = 6349.643
$ws.Range("I62").Value = 4209.8887
$ws.Range("K62").Value = 4209.8887
$ws.Range("M62").Value = -3585.8887

# Row 64
$ws.Range("H64").Value = 10199.333
$ws.Range("J64").Value = 11601
$ws.Range("L64").Value = 11601
$ws.Range("N64").Value = -12097

# Row 65
$ws.Range("H65").Value = 6349.643
$ws.Range("I65").Value = 4209.8887
$ws.Range("K65").Value = 21049.4435
$ws.Range("M65").Value = -17929.4435

# Row 67
$ws.Range("H67").Value = 10199.333
$ws.Range("J67").Value = 11601
$ws.Range("L67").Value = 11601
$ws.Range("N67").Value = -13317

# Row 112
$ws.Range("H112").Value = 1338.8846
$ws.Range("J112").Value = 1374.2273
$ws.Range("L112").Value = 4122.6819
$ws.Range("N112").Value = -6338.6819

# Row 116
$ws.Range("H116").Value = 3944.625
$ws.Range("I116").Value = 3342.6875
$ws.Range("J116").Value = 5148.5
$ws.Range("K116").Value = 3342.6875
$ws.Range("L116").Value = 5148.5
$ws.Range("M116").Value = 99.3125
$ws.Range("N116").Value = -12032.5

# Row 137
$ws.Range("H137").Value = 3059.544
$ws.Range("I137").Value = 2145.5386
$ws.Range("J137").Value = 3329.5908
$ws.Range("K137").Value = 6436.6158
$ws.Range("L137").Value = 9988.7724
$ws.Range("M137").Value = -3886.6158
$ws.Range("N137").Value = -15088.7724

# Row 138
$ws.Range("H138").Value = 2709.8462
$ws.Range("I138").Value = 1267.3103
$ws.Range("J138").Value = 3563.5918
$ws.Range("K138").Value = 3801.9309
$ws.Range("L138").Value = 10690.7754
$ws.Range("M138").Value = 1338.0691
$ws.Range("N138").Value = -20970.7754

# Row 141
$ws.Range("H141").Value = 3063.5715
$ws.Range("I141").Value = 3063.5715
$ws.Range("K141").Value = 9190.7145
$ws.Range("M141").Value = -4010.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3239.8132
$ws.Range("I32").Value = 2157.1
$ws.Range("K32").Value = 2157.1
$ws.Range("M32").Value = -1870.1

# Row 61
$ws.Range("H61").Value = 4359.625
$ws.Range("I61").Value = 3122.074
$ws.Range("J61").Value = 11042.4
$ws.Range("K61").Value = 3122.074
$ws.Range("L61").Value = 11042.4
$ws.Range("M61").Value = -2910.074
$ws.Range("N61").Value = -11466.4

# Row 132
$ws.Range("H132").Value = 1955.3673
$ws.Range("I132").Value = 1540.6666
$ws.Range("K132").Value = 4621.9998
$ws.Range("M132").Value = -2091.9998

# Row 136
$ws.Range("H136").Value = 4359.625
$ws.Range("I136").Value = 3122.074
$ws.Range("J136").Value = 11042.4
$ws.Range("K136").Value = 9366.222
$ws.Range("L136").Value = 33127.2
$ws.Range("M136").Value = -6816.222
$ws.Range("N136").Value = -38227.2

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3776.3157
$ws.Range("I20").Value = 3555.375
$ws.Range("J20").Value = 4954.6665
$ws.Range("K20").Value = 3555.375
$ws.Range("L20").Value = 4954.6665
$ws.Range("M20").Value = -3308.375
$ws.Range("N20").Value = -5448.6665

# Row 99
$ws.Range("H99").Value = 2902
$ws.Range("I99").Value = 2767.9
$ws.Range("J99").Value = 3237.25
$ws.Range("K99").Value = 2767.9
$ws.Range("L99").Value = 3237.25
$ws.Range("M99").Value = -1269.9
$ws.Range("N99").Value = -6233.25

# Row 107
$ws.Range("H107").Value = 2412.2666
$ws.Range("I107").Value = 2207
$ws.Range("K107").Value = 2207
$ws.Range("M107").Value = -287

# Row 109
$ws.Range("H109").Value = 76250.5
$ws.Range("J109").Value = 76250.5
$ws.Range("L109").Value = 76250.5
$ws.Range("N109").Value = -79024.5

# Row 134
$ws.Range("H134").Value = 3155.4666
$ws.Range("I134").Value = 3155.4666
$ws.Range("K134").Value = 9466.399800000001
$ws.Range("M134").Value = -6931.399800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2899.5
$ws.Range("J22").Value = 3739.4
$ws.Range("L22").Value = 3739.4
$ws.Range("N22").Value = -4439.4

# Row 99
$ws.Range("H99").Value = 3983.3333

# Row 126
$ws.Range("H126").Value = 3983.3333

# Row 132
$ws.Range("H132").Value = 3511.7693
$ws.Range("I132").Value = 2216.5
$ws.Range("J132").Value = 5584.2
$ws.Range("K132").Value = 6649.5
$ws.Range("L132").Value = 16752.6
$ws.Range("M132").Value = -4119.5
$ws.Range("N132").Value = -21812.6

# Row 134
$ws.Range("H134").Value = 3163.875
$ws.Range("I134").Value = 2217.5833
$ws.Range("J134").Value = 6002.75
$ws.Range("K134").Value = 6652.749899999999
$ws.Range("L134").Value = 18008.25
$ws.Range("M134").Value = -4117.749899999999
$ws.Range("N134").Value = -23078.25

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 176.6875
$ws.Range("I2").Value = 145.1
$ws.Range("J2").Value = 229.33333
$ws.Range("K2").Value = 870.5999999999999
$ws.Range("L2").Value = 1375.99998
$ws.Range("M2").Value = -757.5999999999999
$ws.Range("N2").Value = -1601.99998

# Row 4
$ws.Range("H4").Value = 7409593.5
$ws.Range("I4").Value = 5882603
$ws.Range("J4").Value = 9572831
$ws.Range("K4").Value = 17647809
$ws.Range("L4").Value = 28718493
$ws.Range("M4").Value = -17647697
$ws.Range("N4").Value = -28718717

# Row 40
$ws.Range("H40").Value = 86.46666999999999
$ws.Range("I40").Value = 17.75
$ws.Range("J40").Value = 165
$ws.Range("K40").Value = 71
$ws.Range("L40").Value = 660
$ws.Range("M40").Value = -2
$ws.Range("N40").Value = -798

# Row 56
$ws.Range("H56").Value = 5359.8
$ws.Range("I56").Value = 5359.8
$ws.Range("K56").Value = 5359.8
$ws.Range("M56").Value = -4829.8

# Row 106
$ws.Range("H106").Value = 14799.6
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 14799.6
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 44398.8
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -46290.8

# Row 107
$ws.Range("H107").Value = 1396.2667
$ws.Range("I107").Value = 671.8333
$ws.Range("J107").Value = 1879.2222
$ws.Range("K107").Value = 2015.4999
$ws.Range("L107").Value = 5637.6666
$ws.Range("M107").Value = -95.49990000000003
$ws.Range("N107").Value = -9477.6666

# Row 111
$ws.Range("H111").Value = 2802.5
$ws.Range("I111").Value = 3157
$ws.Range("J111").Value = 1030
$ws.Range("K111").Value = 9471
$ws.Range("L111").Value = 3090
$ws.Range("M111").Value = -6404
$ws.Range("N111").Value = -9224

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 220723
$ws.Range("J80").Value = 4399.8184
$ws.Range("L80").Value = 4399.8184
$ws.Range("N80").Value = -6395.8184

# Row 83
$ws.Range("H83").Value = 220723
$ws.Range("J83").Value = 4399.8184
$ws.Range("L83").Value = 21999.092
$ws.Range("N83").Value = -31983.092

# Row 102
$ws.Range("H102").Value = 2871.9285
$ws.Range("I102").Value = 2020.8334
$ws.Range("K102").Value = 2020.8334
$ws.Range("M102").Value = -398.8334

# Row 126
$ws.Range("H126").Value = 4557.316
$ws.Range("I126").Value = 3188.6365
$ws.Range("K126").Value = 9565.9095
$ws.Range("M126").Value = -7095.9095

# Row 128
$ws.Range("H128").Value = 70780
$ws.Range("J128").Value = 70780
$ws.Range("L128").Value = 70780
$ws.Range("N128").Value = -80740

# Row 132
$ws.Range("H132").Value = 4034.5
$ws.Range("I132").Value = 1615.5
$ws.Range("K132").Value = 4846.5
$ws.Range("M132").Value = -2316.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5332.222
$ws.Range("I68").Value = 6999
$ws.Range("J68").Value = 4856
$ws.Range("K68").Value = 6999
$ws.Range("L68").Value = 4856
$ws.Range("M68").Value = -6250
$ws.Range("N68").Value = -6354

# Row 71
$ws.Range("H71").Value = 5332.222
$ws.Range("I71").Value = 6999
$ws.Range("J71").Value = 4856
$ws.Range("K71").Value = 34995
$ws.Range("L71").Value = 24280
$ws.Range("M71").Value = -31251
$ws.Range("N71").Value = -31768

# Row 122
$ws.Range("H122").Value = 273961.12
$ws.Range("J122").Value = 9802.375
$ws.Range("L122").Value = 29407.125
$ws.Range("N122").Value = -34307.125

# Row 127
$ws.Range("H127").Value = 72953.336
$ws.Range("J127").Value = 72953.336
$ws.Range("L127").Value = 72953.336
$ws.Range("N127").Value = -82873.336

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2563.3333
$ws.Range("I122").Value = 1780.8182
$ws.Range("J122").Value = 3101.3125
$ws.Range("K122").Value = 5342.4546
$ws.Range("L122").Value = 9303.9375
$ws.Range("M122").Value = -2892.4546
$ws.Range("N122").Value = -14203.9375

# Row 132
$ws.Range("H132").Value = 4726.5
$ws.Range("I132").Value = 1870.8
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 4856.5
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -3082.4
